$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.34 = 8807.99 pesos`n✅ 8807.99 pesos = 2.34 = 958.41 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsTasas.Range("N10").Value = 427
$wsTasas.Range("O10").Value = 3761.01
$wsTasas.Range("N12").Value = 3770
$wsTasas.Range("O12").Value = 410.22
